$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "publish" column (D) currently stores dates as plain text strings
# ("2018-03-06", "2018-31-12"). Convert them to real Excel date values,
# formatted with a standard short-date number format, while preserving
# the existing visual styling (border/wrap on the header-adjacent row,
# vertical-centered alignment, etc.) of each cell.

$d2 = Get-Date -Year 2018 -Month 3 -Day 6 -Hour 0 -Minute 0 -Second 0
$d3 = Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0

$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Value = $d2

$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("D3").Value = $d3

# Move the active selection to D1 (was previously left on E12).
$ws.Range("D1").Select()
